$wb = $excel.ActiveWorkbook

# Add the new "DAC101C081" worksheet after the last sheet (TMF8801) so it
# lands at the end of the tab strip and becomes the active/selected sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "DAC101C081"

# Header row (shared with the other register sheets).
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Hex Address"
$ws.Range("C1").Value = "Default Value"
$ws.Range("D1").Value = "Bit Width"
$ws.Range("E1").Value = "Bit Index (High)"
$ws.Range("F1").Value = "Bit Index (Low)"

# Row 2 - PD field.
$ws.Range("A2").Value = "PD"
$ws.Range("B2").Value = "None"
$ws.Range("C2").Value = "0x0000"
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = 12

# Row 3 - DATA field.
$ws.Range("A3").Value = "DATA"
$ws.Range("B3").Value = "None"
$ws.Range("C3").Value = "0x0000"
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = 11
$ws.Range("F3").Value = 2

# The "Default Value" column is tagged with the workbook's "Bad" cell style
# (as used elsewhere for a "None" register address) but with the red
# text/pink fill cleared back to the normal look.
$ws.Range("C2").Style = "Bad"
$ws.Range("C2").Font.ThemeColor = 1
$ws.Range("C2").Interior.Pattern = -4142

$ws.Range("C3").Style = "Bad"
$ws.Range("C3").Font.ThemeColor = 1
$ws.Range("C3").Interior.Pattern = -4142

$ws.Columns.Item(1).ColumnWidth = 20.3

$ws.Range("B3").Select() | Out-Null
